$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Protect column D (Price) cells from Excel auto-converting numeric-looking
# text (e.g. "55.183.15", "139.43") into actual numbers / floats.
# We force Text format before assigning, then restore the default "Normal"
# style afterwards so no stray style/number-format remains on the cells.
$dRange = $ws.Range("D2:D51")
$dRange.NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = "55.183.15"
$ws.Range("E2").Value = "  +2.55%  "

# Row 3
$ws.Range("D3").Value = "2.488.40"
$ws.Range("E3").Value = "  +6.75%  "

# Row 4
$ws.Range("E4").Value = "  +0.29%  "

# Row 5
$ws.Range("D5").Value = "478.88"
$ws.Range("E5").Value = "  +8.10%  "

# Row 6
$ws.Range("D6").Value = "139.43"
$ws.Range("E6").Value = "  +10.81%  "

# Row 7
$ws.Range("E7").Value = "  +0.36%  "

# Row 8
$ws.Range("D8").Value = "0.507"
$ws.Range("E8").Value = "  +6.14%  "

# Row 9
$ws.Range("D9").Value = "2.487.94"
$ws.Range("E9").Value = "  +6.59%  "

# Row 10
$ws.Range("D10").Value = "0.0978"
$ws.Range("E10").Value = "  +6.36%  "

# Row 11
$ws.Range("D11").Value = "5.44"
$ws.Range("E11").Value = "  +1.42%  "

# Row 12
$ws.Range("D12").Value = "0.326"
$ws.Range("E12").Value = "  +5.45%  "

# Row 13
$ws.Range("E13").Value = "  +0.31%  "

# Row 14
$ws.Range("D14").Value = "2.925.47"
$ws.Range("E14").Value = "  +7.59%  "

# Row 15
$ws.Range("D15").Value = "55.299.35"
$ws.Range("E15").Value = "  +2.76%  "

# Row 16
$ws.Range("D16").Value = "20.47"
$ws.Range("E16").Value = "  +8.37%  "

# Row 17
$ws.Range("D17").Value = "0.0000136"
$ws.Range("E17").Value = "  +12.29%  "

# Row 18
$ws.Range("D18").Value = "2.495.25"
$ws.Range("E18").Value = "  +6.07%  "

# Row 19
$ws.Range("D19").Value = "4.36"
$ws.Range("E19").Value = "  +10.63%  "

# Row 20
$ws.Range("D20").Value = "318.69"
$ws.Range("E20").Value = "  +6.66%  "

# Row 21
$ws.Range("D21").Value = "9.98"
$ws.Range("E21").Value = "  +8.66%  "

# Row 22
$ws.Range("D22").Value = "0.996"
$ws.Range("E22").Value = "  -0.19%  "

# Row 23
$ws.Range("D23").Value = "5.66"
$ws.Range("E23").Value = "  +5.02%  "

# Row 24
$ws.Range("D24").Value = "57.69"
$ws.Range("E24").Value = "  +3.48%  "

# Row 25
$ws.Range("B25").Value = "Kaspa"
$ws.Range("C25").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D25").Value = "0.166"
$ws.Range("E25").Value = "  +8.82%  "

# Row 26
$ws.Range("B26").Value = "Polygon"
$ws.Range("C26").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D26").Value = "0.406"
$ws.Range("E26").Value = "  +10.55%  "

# Row 27
$ws.Range("D27").Value = "1.01"
$ws.Range("E27").Value = "  +0.94%  "

# Row 28
$ws.Range("D28").Value = "2.605.28"
$ws.Range("E28").Value = "  +6.99%  "

# Row 29
$ws.Range("D29").Value = "7.36"
$ws.Range("E29").Value = "  +5.37%  "

# Row 30
$ws.Range("D30").Value = "0.0₃0778"
$ws.Range("E30").Value = "  +10.61%  "

# Row 31
$ws.Range("E31").Value = "  +0.39%  "

# Row 32
$ws.Range("D32").Value = "148.74"
$ws.Range("E32").Value = "  +1.54%  "

# Row 33
$ws.Range("D33").Value = "18.12"
$ws.Range("E33").Value = "  +5.01%  "

# Row 34
$ws.Range("D34").Value = "1.47"
$ws.Range("E34").Value = "  +8.73%  "

# Row 35
$ws.Range("D35").Value = "5.17"
$ws.Range("E35").Value = "  +10.50%  "

# Row 36
$ws.Range("D36").Value = "3.68"
$ws.Range("E36").Value = "  +3.42%  "

# Row 37
$ws.Range("D37").Value = "1.11"
$ws.Range("E37").Value = "  +9.87%  "

# Row 38
$ws.Range("D38").Value = "0.846"
$ws.Range("E38").Value = "  +1.65%  "

# Row 39
$ws.Range("D39").Value = "34.22"
$ws.Range("E39").Value = "  +2.76%  "

# Row 40
$ws.Range("D40").Value = "0.997"
$ws.Range("E40").Value = "  +0.29%  "

# Row 41
$ws.Range("D41").Value = "0.608"
$ws.Range("E41").Value = "  +16.32%  "

# Row 42
$ws.Range("D42").Value = "0.0550"
$ws.Range("E42").Value = "  +10.65%  "

# Row 43
$ws.Range("D43").Value = "3.38"
$ws.Range("E43").Value = "  +7.49%  "

# Row 44
$ws.Range("D44").Value = "1.31"
$ws.Range("E44").Value = "  +8.45%  "

# Row 45
$ws.Range("D45").Value = "10.16"
$ws.Range("E45").Value = "  -1.19%  "

# Row 46
$ws.Range("D46").Value = "1.972.26"
$ws.Range("E46").Value = "  +1.38%  "

# Row 47
$ws.Range("D47").Value = "0.0901"
$ws.Range("E47").Value = "  +8.37%  "

# Row 48
$ws.Range("D48").Value = "0.0222"
$ws.Range("E48").Value = "  +5.85%  "

# Row 49
$ws.Range("D49").Value = "246.50"
$ws.Range("E49").Value = "  +28.59%  "

# Row 50
$ws.Range("D50").Value = "4.48"
$ws.Range("E50").Value = "  +12.89%  "

# Row 51
$ws.Range("D51").Value = "17.39"
$ws.Range("E51").Value = "  +8.94%  "

# Restore default styling on column D so number formats match the original file
$dRange.Style = "Normal"
